$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 244 down to the new rows 245-247
$ws.Range("A244:D244").Copy()
$ws.Range("A245:D247").PasteSpecial(-4122)

$ws.Range("A245").Value = 44319
$ws.Range("B245").Value = 0
$ws.Range("C245").Value = 2
$ws.Range("D245").Value = 87.56567425569177

$ws.Range("A246").Value = 44320
$ws.Range("B246").Value = 2
$ws.Range("C246").Value = 3
$ws.Range("D246").Value = 131.3485113835376

$ws.Range("A247").Value = 44321
$ws.Range("B247").Value = 0
$ws.Range("C247").Value = 3
$ws.Range("D247").Value = 131.3485113835376
